$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loans")

# Row 4 picked up an explicit (custom) row height in the new revision.
$ws.Rows("4").RowHeight = 14.25

# Update row 5's loan dates to the new timestamps, then propagate the
# whole row (A:E) onto rows 6 and 7 via Copy so the string cells keep
# their shared-string typing (no stray numeric/style coercion).
$ws.Range("D5").Value = "Wed May 23 13:38:51 GMT-03:00 2018"
$ws.Range("E5").Value = "Thu May 24 13:38:51 GMT-03:00 2018"
$ws.Range("A5:E5").Copy($ws.Range("A6:E6"))
$ws.Range("A5:E5").Copy($ws.Range("A7:E7"))

# The trailing loan rows (8-11) were removed from the sheet.
$ws.Rows("8:11").Delete()
